# edit.ps1 - "update RE ADME file"
#
# Adds a title-page header block (name / id + spacer paragraphs) in front
# of the existing content, and restores a set of "ListLabel 10".."ListLabel 18"
# character styles (with a couple of docDefaults tweaks) to the style sheet,
# mirroring what LibreOffice produces when new numbered lists are added to a
# document that already defines ListLabel 1..9.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert the five new paragraphs at the very top of the document body.
#    Using Range(0,0).InsertXML lets us express the exact pPr/rPr markup
#    (alignment, bold/underline/size run formatting, empty marker runs)
#    for each paragraph instead of relying on a single blended style.
# ---------------------------------------------------------------------
$introXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="center"/><w:rPr><w:b/><w:b/><w:bCs/><w:sz w:val="36"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:u w:val="single"/></w:rPr><w:t>Name: Issa Negasa</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:u w:val="single"/></w:rPr><w:t>Id: 984556</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="left"/><w:rPr><w:b/><w:b/><w:bCs/><w:sz w:val="36"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="left"/><w:rPr><w:b/><w:b/><w:bCs/><w:sz w:val="36"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$topRange = $d.Range(0, 0)
$topRange.InsertXML($introXml)

# ---------------------------------------------------------------------
# 2) Recreate the missing "ListLabel 10".."ListLabel 18" character styles
#    (the document already has ListLabel 1..9 with this same shape: a
#    qFormat character style whose only explicit run formatting is the
#    OpenSymbol complex-script font used for bullet glyphs).
# ---------------------------------------------------------------------
for ($i = 10; $i -le 18; $i++) {
    $style = $d.Styles.Add("ListLabel " + $i, 2)
    $style.Font.NameBi = "OpenSymbol"
    $style.QuickStyle = 1
}

Write-Output "Inserted header paragraphs and added ListLabel10-18 styles."
